$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column N (14th column) to make room for the
# new "PAN" / "{vendor:pan_no}" field, shifting existing columns N.. to the right.
$ws.Columns.Item(14).Insert()

# Populate the new column's header (row1) and sample/placeholder (row2) values.
$ws.Cells.Item(1, 14).Value = "PAN"
$ws.Cells.Item(2, 14).Value = "{vendor:pan_no}"
